$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.247.52'
$ws.Range("E2").Value = '  +0.72%  '

$ws.Range("D3").Value = '2.226.44'
$ws.Range("E3").Value = '  +0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.45'
$ws.Range("E5").Value = '  +0.10%  '

$ws.Range("E6").Value = '  +0.56%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.98'
$ws.Range("E7").Value = '  +0.46%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.614'
$ws.Range("E9").Value = '  +0.63%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.43'
$ws.Range("E10").Value = '  +5.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0976'
$ws.Range("E11").Value = '  +4.55%  '

$ws.Range("E12").Value = '  +1.78%  '

$ws.Range("E13").Value = '  +1.35%  '

$ws.Range("E14").Value = '  +0.19%  '

$ws.Range("E15").Value = '  +1.25%  '

$ws.Range("D16").Value = '2.207.29'
$ws.Range("E16").Value = '  -1.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000116'
$ws.Range("E17").Value = '  +20.72%  '

$ws.Range("D18").Value = '42.114.68'
$ws.Range("E18").Value = '  +0.62%  '

$ws.Range("E19").Value = '  +2.45%  '

$ws.Range("E20").Value = '  +1.68%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.99'
$ws.Range("E21").Value = '  +39.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.14'
$ws.Range("E22").Value = '  +1.27%  '

$ws.Range("E23").Value = '  -3.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.87'
$ws.Range("E24").Value = '  +8.51%  '

$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("E26").Value = '  -1.10%  '

$ws.Range("E27").Value = '  +2.04%  '

$ws.Range("E28").Value = '  +3.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.85'
$ws.Range("E29").Value = '  -0.38%  '

$ws.Range("E30").Value = '  +3.70%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.69'
$ws.Range("E31").Value = '  +19.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0801'
$ws.Range("E32").Value = '  -0.64%  '

$ws.Range("E33").Value = '  +1.20%  '

$ws.Range("E34").Value = '  +1.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '29.25'
$ws.Range("E35").Value = '  -2.68%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.40'
$ws.Range("E36").Value = '  +1.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0303'
$ws.Range("E37").Value = '  +3.33%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '13.01'
$ws.Range("E38").Value = '  -1.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.62'
$ws.Range("E40").Value = '  -1.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '62.52'
$ws.Range("E41").Value = '  +5.54%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.200'
$ws.Range("E42").Value = '  +1.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.80'
$ws.Range("E43").Value = '  +2.54%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '105.35'
$ws.Range("E44").Value = '  -4.48%  '

$ws.Range("E45").Value = '  +3.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.993'
$ws.Range("E46").Value = '  -0.11%  '

$ws.Range("E47").Value = '  +7.73%  '

$ws.Range("E48").Value = '  +1.82%  '

$ws.Range("E49").Value = '  +2.78%  '

$ws.Range("E50").Value = '  +0.74%  '

$ws.Range("E51").Value = '  +1.05%  '
